# Insert a new "AnatomicalEntity" sheet right before "Organization"
# (i.e. right after "NamedThing"), with the same header row / schema
# that is used by the other top-level NamedThing subclasses.

$wb = $excel.ActiveWorkbook

$organization = $wb.Worksheets.Item("Organization")
$newSheet = $wb.Worksheets.Add($organization)
$newSheet.Name = "AnatomicalEntity"

$headers = @(
    "id",
    "category",
    "name",
    "description",
    "subclass_of",
    "related_to",
    "contributor_name",
    "contributor_github_name",
    "contributor_orcid",
    "contribution_date"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
